# "Changes on Names For Paper"
# Rename the Tech values in column B to the new naming convention and
# remove the centered style from those cells (matches the canonical diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "Transmission_Interregional" = "TRANSMISSION_INTERREGIONAL"
    "OffshoreWind_New_C6"        = "WIND-OFFSHORE-C6_NEW"
    "LandWind_New_C8"            = "WIND-LAND-C8_NEW"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
    $cell.Style = "Normal"
}

# Move the active selection as recorded in the saved workbook.
$ws.Range("B9").Select()
